# Auto-generated edit script applying the Kujata_Profits market-data refresh diff
# Updates per-leve market price/profit figures across all crafting class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16627.375
$ws.Range("I21").Value = 19250
$ws.Range("J21").Value = 14004.75
$ws.Range("K21").Value = 19250
$ws.Range("L21").Value = 14004.75
$ws.Range("M21").Value = -18782
$ws.Range("N21").Value = -14940.75
$ws.Range("H23").Value = 16627.375
$ws.Range("I23").Value = 19250
$ws.Range("J23").Value = 14004.75
$ws.Range("K23").Value = 19250
$ws.Range("L23").Value = 14004.75
$ws.Range("M23").Value = -19016
$ws.Range("N23").Value = -14472.75
$ws.Range("H112").Value = 2406.0476
$ws.Range("J112").Value = 2807.647
$ws.Range("L112").Value = 8422.940999999999
$ws.Range("N112").Value = -10638.941
$ws.Range("H129").Value = 740.8182
$ws.Range("I129").Value = 375.2857
$ws.Range("J129").Value = 911.4
$ws.Range("K129").Value = 1125.8571
$ws.Range("L129").Value = 2734.2
$ws.Range("M129").Value = 3874.1429
$ws.Range("N129").Value = -12734.2
$ws.Range("H138").Value = 1788.04
$ws.Range("I138").Value = 988.95654
$ws.Range("J138").Value = 2026.7273
$ws.Range("K138").Value = 2966.86962
$ws.Range("L138").Value = 6080.1819
$ws.Range("M138").Value = 2173.13038
$ws.Range("N138").Value = -16360.1819
$ws.Range("H140").Value = 35914
$ws.Range("J140").Value = 35914
$ws.Range("L140").Value = 35914
$ws.Range("N140").Value = -46274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1235.65
$ws.Range("I2").Value = 994.0769
$ws.Range("J2").Value = 1684.2858
$ws.Range("K2").Value = 994.0769
$ws.Range("L2").Value = 1684.2858
$ws.Range("M2").Value = -881.0769
$ws.Range("N2").Value = -1910.2858
$ws.Range("H4").Value = 722.8
$ws.Range("I4").Value = 450
$ws.Range("K4").Value = 450
$ws.Range("M4").Value = -334
$ws.Range("H16").Value = 1700
$ws.Range("I16").Value = 1700
$ws.Range("K16").Value = 1700
$ws.Range("M16").Value = -1413
$ws.Range("H32").Value = 8675.552
$ws.Range("I32").Value = 6753.9595
$ws.Range("J32").Value = 19613.846
$ws.Range("K32").Value = 6753.9595
$ws.Range("L32").Value = 19613.846
$ws.Range("M32").Value = -6466.9595
$ws.Range("N32").Value = -20187.846
$ws.Range("H38").Value = 8354
$ws.Range("I38").Value = 5020
$ws.Range("J38").Value = 10021
$ws.Range("K38").Value = 5020
$ws.Range("L38").Value = 10021
$ws.Range("M38").Value = -4553
$ws.Range("N38").Value = -10955
$ws.Range("H61").Value = 41668124
$ws.Range("I61").Value = 52632620
$ws.Range("K61").Value = 52632620
$ws.Range("M61").Value = -52632408
$ws.Range("H74").Value = 2924
$ws.Range("I74").Value = 2550.375
$ws.Range("J74").Value = 3137.5
$ws.Range("K74").Value = 2550.375
$ws.Range("L74").Value = 3137.5
$ws.Range("M74").Value = -1676.375
$ws.Range("N74").Value = -4885.5
$ws.Range("H77").Value = 2924
$ws.Range("I77").Value = 2550.375
$ws.Range("J77").Value = 3137.5
$ws.Range("K77").Value = 12751.875
$ws.Range("L77").Value = 15687.5
$ws.Range("M77").Value = -8383.875
$ws.Range("N77").Value = -24423.5
$ws.Range("H116").Value = 1235.65
$ws.Range("I116").Value = 994.0769
$ws.Range("J116").Value = 1684.2858
$ws.Range("K116").Value = 994.0769
$ws.Range("L116").Value = 1684.2858
$ws.Range("M116").Value = 1299.9231
$ws.Range("N116").Value = -6272.2858
$ws.Range("H122").Value = 2285.3684
$ws.Range("I122").Value = 1970.2354
$ws.Range("K122").Value = 5910.706200000001
$ws.Range("M122").Value = -3460.706200000001
$ws.Range("H132").Value = 2459.2708
$ws.Range("I132").Value = 1732.6666
$ws.Range("K132").Value = 5197.9998
$ws.Range("M132").Value = -2667.9998
$ws.Range("H136").Value = 41668124
$ws.Range("I136").Value = 52632620
$ws.Range("K136").Value = 157897860
$ws.Range("M136").Value = -157895310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1235.65
$ws.Range("I3").Value = 994.0769
$ws.Range("J3").Value = 1684.2858
$ws.Range("K3").Value = 994.0769
$ws.Range("L3").Value = 1684.2858
$ws.Range("M3").Value = -880.0769
$ws.Range("N3").Value = -1912.2858
$ws.Range("H19").Value = 4999.3335
$ws.Range("I19").Value = 4999.3335
$ws.Range("K19").Value = 4999.3335
$ws.Range("M19").Value = -4826.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1152.4657
$ws.Range("I31").Value = 1100.7705
$ws.Range("J31").Value = 1415.25
$ws.Range("K31").Value = 1100.7705
$ws.Range("L31").Value = 1415.25
$ws.Range("M31").Value = -805.7705000000001
$ws.Range("N31").Value = -2005.25
$ws.Range("H34").Value = 1152.4657
$ws.Range("I34").Value = 1100.7705
$ws.Range("J34").Value = 1415.25
$ws.Range("K34").Value = 1100.7705
$ws.Range("L34").Value = 1415.25
$ws.Range("M34").Value = -898.7705000000001
$ws.Range("N34").Value = -1819.25
$ws.Range("H132").Value = 2166.182
$ws.Range("I132").Value = 1572
$ws.Range("J132").Value = 2879.2
$ws.Range("K132").Value = 4716
$ws.Range("L132").Value = 8637.599999999999
$ws.Range("M132").Value = -2186
$ws.Range("N132").Value = -13697.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 405.83334
$ws.Range("I7").Value = 463.46155
$ws.Range("J7").Value = 256
$ws.Range("K7").Value = 1390.38465
$ws.Range("L7").Value = 768
$ws.Range("M7").Value = -1278.38465
$ws.Range("N7").Value = -992
$ws.Range("H59").Value = 3499.9
$ws.Range("I59").Value = 1499.5
$ws.Range("J59").Value = 4000
$ws.Range("K59").Value = 4498.5
$ws.Range("L59").Value = 12000
$ws.Range("M59").Value = -3958.5
$ws.Range("N59").Value = -13080
$ws.Range("H113").Value = 651.2558
$ws.Range("I113").Value = 576.6923
$ws.Range("J113").Value = 683.56665
$ws.Range("K113").Value = 1730.0769
$ws.Range("L113").Value = 2050.69995
$ws.Range("M113").Value = 439.9231
$ws.Range("N113").Value = -6390.69995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 860.5
$ws.Range("I31").Value = 860.5
$ws.Range("K31").Value = 860.5
$ws.Range("M31").Value = -568.5
$ws.Range("H37").Value = 860.5
$ws.Range("I37").Value = 860.5
$ws.Range("K37").Value = 860.5
$ws.Range("M37").Value = -583.5
$ws.Range("H132").Value = 2526.75
$ws.Range("I132").Value = 2205.4666
$ws.Range("K132").Value = 6616.399800000001
$ws.Range("M132").Value = -4086.399800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 7225
$ws.Range("I32").Value = 6966.6665
$ws.Range("K32").Value = 6966.6665
$ws.Range("M32").Value = -6649.6665
$ws.Range("H68").Value = 1282.1428
$ws.Range("I68").Value = 1273.0769
$ws.Range("J68").Value = 1400
$ws.Range("K68").Value = 1273.0769
$ws.Range("L68").Value = 1400
$ws.Range("M68").Value = -524.0769
$ws.Range("N68").Value = -2898
$ws.Range("H71").Value = 1282.1428
$ws.Range("I71").Value = 1273.0769
$ws.Range("J71").Value = 1400
$ws.Range("K71").Value = 6365.3845
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = -2621.3845
$ws.Range("N71").Value = -14488
$ws.Range("H122").Value = 9616822
$ws.Range("I122").Value = 12501373
$ws.Range("K122").Value = 37504119
$ws.Range("M122").Value = -37501669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1500
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H54").Value = 25000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H122").Value = 9260505
$ws.Range("I122").Value = 11906023
$ws.Range("J122").Value = 1193.3334
$ws.Range("K122").Value = 35718069
$ws.Range("L122").Value = 3580.0002
$ws.Range("M122").Value = -35715619
$ws.Range("N122").Value = -8480.0002
$ws.Range("H141").Value = 62500
$ws.Range("J141").Value = 62500
$ws.Range("L141").Value = 62500
$ws.Range("N141").Value = -72860
